$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Diary entries for the past two weeks (rows 22-25), following the
# same layout/formatting used by the existing entries (rows 10-21).
# ------------------------------------------------------------------

# Seed formatting for the four new data rows from an existing plain
# data row (row 10: Date | Time-as-text | Participants | Goal |
# Achievements | Reflection | Mood) so borders/fills/fonts line up.
$ws.Range("A10:G10").Copy()
$ws.Range("A22:G25").PasteSpecial(-4122)

# Seed formatting for the trailing blank rows from the existing blank
# row (row 22 originally, now shifted — use row 21 as a stand-in data
# row and then blank the values) so the blank rows keep the table's
# border/fill styling.
$ws.Range("A21:G21").Copy()
$ws.Range("A26:G31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Row 22 : 2/12/2020 ----
$ws.Range("A22").Value = 43873
$ws.Range("A22").NumberFormat = "m/d/yyyy"
$ws.Range("B22").Value = "9:00PM - 1:00AM"
$ws.Range("C22").Value = "N/A"
$ws.Range("D22").Value = "Study for the midterm"
$ws.Range("E22").Value = "Hopefully memorized everything"
$ws.Range("F22").Value = "There weren't a ton of materials so I should be fine"
$ws.Range("G22").Value = "Decent, tired"

# ---- Row 23 : 2/13/2020 ----
$ws.Range("A23").Value = 43874
$ws.Range("A23").NumberFormat = "m/d/yyyy"
$ws.Range("B23").Value = "5:00 PM - 7:50PM"
$ws.Range("C23").Value = "N/A"
$ws.Range("D23").Value = "Take midterm, and learn"
$ws.Range("E23").Value = "Did okay on the midterm I think, Learned about the stakeholders, functionalities, and key developers of a project"
$ws.Range("F23").Value = "Finding all these should be easy for Runelite, since I am very familiar with the game and the client"
$ws.Range("G23").Value = "Exhausted from the midterm"

# ---- Row 24 : 2/19/2020 ----
$ws.Range("A24").Value = 43880
$ws.Range("A24").NumberFormat = "m/d/yyyy"
$ws.Range("B24").Value = "9:00PM - 2:00AM"
$ws.Range("C24").Value = "Thuc, Harry"
$ws.Range("D24").Value = "Work on the Big Picture"
$ws.Range("E24").Value = "Was able to identify everything and explain them all"
$ws.Range("F24").Value = "Instructions were clear, but not sure if bullet point answers would suffice, so we tried to explain everything in detail"
$ws.Range("G24").Value = "Confused?"

# ---- Row 25 : 2/20/2020 ----
$ws.Range("A25").Value = 43881
$ws.Range("A25").NumberFormat = "m/d/yyyy"
$ws.Range("B25").Value = "3:00PM - 5:00PM"
$ws.Range("C25").Value = "Thuc, Harry"
$ws.Range("D25").Value = "Finishing touches to the homework"
$ws.Range("E25").Value = "pull request submitted"
$ws.Range("F25").Value = "Nothing much to reflect on, hopefully meets expectations"
$ws.Range("G25").Value = "Uncertain?"

# ------------------------------------------------------------------
# Formatting touch-ups to match how the sheet actually ended up:
# ------------------------------------------------------------------

# Column B (Time) for the new rows is free text like "9:00PM - 1:00AM"
# rather than a time-of-day value, so it keeps the italic green style
# but right-aligned, and no time number format.
$ws.Range("B22:B31").NumberFormat = "General"
$ws.Range("B22:B31").HorizontalAlignment = -4152
$ws.Range("B22:B31").Font.Italic = $true
$ws.Range("B22:B31").Font.Size = 12
$ws.Range("B22:B31").Font.Color = 24832

# G23:G25 use the same 12pt italic green style as the rest of the row
# (rather than the smaller 11pt "mood" style used elsewhere).
$ws.Range("G23:G25").Font.Size = 12

# The trailing blank rows (26-31) pick up a slightly different italic
# green font (12pt, Calibri) together with the date number format.
$ws.Range("A26:G31").Font.Name = "Calibri"
$ws.Range("A26:G31").Font.Size = 12
$ws.Range("A26:G31").Font.Italic = $true
$ws.Range("A26:G31").Font.Color = 24832
$ws.Range("A26:G31").NumberFormat = "m/d/yyyy"

# Column B across the blank rows keeps the plain (non-Calibri) italic
# green style used for the rest of the new Time entries above.
$ws.Range("B26:B31").Font.Name = "Arial"
$ws.Range("B26:B31").NumberFormat = "General"

# C30 is a one-off: it kept the plain (non date-formatted) style.
$ws.Range("C30").Font.Name = "Arial"
$ws.Range("C30").NumberFormat = "General"
